# Auto-generated PowerShell Excel COM-interop edit script
# Applies numeric corrections to existing rows and appends new rows 146-161

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: update existing numeric values (re-calculated statistics) ---
$ws.Range("F12").Value = 4.72355
$ws.Range("G12").Value = 5.79667068965517
$ws.Range("L12").Value = 3.38695
$ws.Range("F13").Value = 4.72355
$ws.Range("G13").Value = 5.79667068965517
$ws.Range("L13").Value = 3.38695
$ws.Range("G24").Value = 4.39684043691337
$ws.Range("G25").Value = 4.39684043691337
$ws.Range("G26").Value = 0.788226562319772
$ws.Range("G27").Value = 0.788226562319772
$ws.Range("G28").Value = 5.54195550679922
$ws.Range("G29").Value = 5.54195550679922
$ws.Range("G40").Value = 4.86334335897353
$ws.Range("G41").Value = 4.86334335897353
$ws.Range("G42").Value = 0.764976562319772
$ws.Range("G43").Value = 0.764976562319772
$ws.Range("G44").Value = 5.02187217346589
$ws.Range("H44").Value = 15.7491
$ws.Range("L44").Value = 2.4798
$ws.Range("M44").Value = 7.70927
$ws.Range("G45").Value = 5.02187217346589
$ws.Range("H45").Value = 15.7491
$ws.Range("L45").Value = 2.4798
$ws.Range("M45").Value = 7.70927
$ws.Range("G56").Value = 4.52352785978244
$ws.Range("G57").Value = 4.52352785978244
$ws.Range("G58").Value = 0.740593228986438
$ws.Range("G59").Value = 0.740593228986438
$ws.Range("F60").Value = 3.24195
$ws.Range("G60").Value = 4.17237217346589
$ws.Range("L60").Value = 2.4798
$ws.Range("M60").Value = 6.63813
$ws.Range("N60").Value = 8.70327
$ws.Range("F61").Value = 3.24195
$ws.Range("G61").Value = 4.17237217346589
$ws.Range("L61").Value = 2.4798
$ws.Range("M61").Value = 6.63813
$ws.Range("N61").Value = 8.70327
$ws.Range("G68").Value = 2537.47217943343
$ws.Range("I68").Value = 14008.66538
$ws.Range("G69").Value = 2537.47217943343
$ws.Range("I69").Value = 14008.66538
$ws.Range("G70").Value = 2537.47217943343
$ws.Range("I70").Value = 14008.66538
$ws.Range("G71").Value = 2537.47217943343
$ws.Range("I71").Value = 14008.66538
$ws.Range("G72").Value = 4.10092571716825
$ws.Range("G73").Value = 4.10092571716825
$ws.Range("G74").Value = 0.785734997756939
$ws.Range("G75").Value = 0.785734997756939
$ws.Range("G76").Value = 3.92424550679922
$ws.Range("L76").Value = 2.35915
$ws.Range("M76").Value = 5.78912
$ws.Range("N76").Value = 8.241630000000001
$ws.Range("G77").Value = 3.92424550679922
$ws.Range("L77").Value = 2.35915
$ws.Range("M77").Value = 5.78912
$ws.Range("N77").Value = 8.241630000000001
$ws.Range("G84").Value = 4534.52217943343
$ws.Range("I84").Value = 25000
$ws.Range("N84").Value = 15205.19923
$ws.Range("G85").Value = 4534.52217943343
$ws.Range("I85").Value = 25000
$ws.Range("N85").Value = 15205.19923
$ws.Range("G86").Value = 4534.52217943343
$ws.Range("I86").Value = 25000
$ws.Range("N86").Value = 15205.19923
$ws.Range("G87").Value = 4534.52217943343
$ws.Range("I87").Value = 25000
$ws.Range("N87").Value = 15205.19923
$ws.Range("G88").Value = 5.61710258533583
$ws.Range("G89").Value = 5.61710258533583
$ws.Range("G90").Value = 0.795393331090273
$ws.Range("G91").Value = 0.795393331090273
$ws.Range("G92").Value = 5.06769050679922
$ws.Range("M92").Value = 8.227270000000001
$ws.Range("G93").Value = 5.06769050679922
$ws.Range("M93").Value = 8.227270000000001
$ws.Range("G100").Value = 5034.58733703194
$ws.Range("I100").Value = 24278.19495
$ws.Range("N100").Value = 15683.81277
$ws.Range("G101").Value = 5034.58733703194
$ws.Range("I101").Value = 24278.19495
$ws.Range("N101").Value = 15683.81277
$ws.Range("G102").Value = 5034.58733703194
$ws.Range("I102").Value = 24278.19495
$ws.Range("N102").Value = 15683.81277
$ws.Range("G103").Value = 5034.58733703194
$ws.Range("I103").Value = 24278.19495
$ws.Range("N103").Value = 15683.81277
$ws.Range("G104").Value = 6.29002981837308
$ws.Range("G105").Value = 6.29002981837308
$ws.Range("G106").Value = 0.870122137732713
$ws.Range("G107").Value = 0.870122137732713
$ws.Range("G116").Value = 5131.02801499804
$ws.Range("I116").Value = 24278.19495
$ws.Range("N116").Value = 15900
$ws.Range("G117").Value = 5131.02801499804
$ws.Range("I117").Value = 24278.19495
$ws.Range("N117").Value = 15900
$ws.Range("G118").Value = 5131.02801499804
$ws.Range("I118").Value = 24278.19495
$ws.Range("N118").Value = 15900
$ws.Range("G119").Value = 5131.02801499804
$ws.Range("I119").Value = 24278.19495
$ws.Range("N119").Value = 15900
$ws.Range("G120").Value = 7.13211587040832
$ws.Range("G121").Value = 7.13211587040832
$ws.Range("G122").Value = 0.903257730953052
$ws.Range("G123").Value = 0.903257730953052
$ws.Range("G132").Value = 5361.6956646471
$ws.Range("I132").Value = 25105.26274
$ws.Range("N132").Value = 12470
$ws.Range("G133").Value = 5361.6956646471
$ws.Range("I133").Value = 25105.26274
$ws.Range("N133").Value = 12470
$ws.Range("G134").Value = 5361.6956646471
$ws.Range("I134").Value = 25105.26274
$ws.Range("N134").Value = 12470
$ws.Range("G135").Value = 5361.6956646471
$ws.Range("I135").Value = 25105.26274
$ws.Range("N135").Value = 12470
$ws.Range("G138").Value = 0.888687826775966
$ws.Range("G139").Value = 0.888687826775966

# --- Part 2: append new rows 146-161 (2019-2023 period) ---
# Row 146
$ws.Range("A146").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B146").Value = "DRP (95th Percentile)"
$ws.Range("C146").Value = "D"
$ws.Range("D146").Value = "2019 - 2023"
$ws.Range("E146").Value = "Impact"
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 1.22275438596491
$ws.Range("H146").Value = 3.28
$ws.Range("I146").Value = 2.852
$ws.Range("L146").Value = 0.881
$ws.Range("M146").Value = 2.1048
$ws.Range("N146").Value = 2.7462
$ws.Range("O146").Value = 1784937.706
$ws.Range("P146").Value = 5565566.687
$ws.Range("Q146").Value = "Rangitikei District"
$ws.Range("R146").Value = "Rangitīkei-Turakina"
$ws.Range("S146").Value = "Turakina"
$ws.Range("T146").Value = "Tura_1c"
$ws.Range("U146").Value = "mg/L"
# Row 147
$ws.Range("A147").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B147").Value = "DRP (Median)"
$ws.Range("C147").Value = "D"
$ws.Range("D147").Value = "2019 - 2023"
$ws.Range("E147").Value = "Impact"
$ws.Range("F147").Value = 1
$ws.Range("G147").Value = 1.22275438596491
$ws.Range("H147").Value = 3.28
$ws.Range("I147").Value = 2.852
$ws.Range("L147").Value = 0.881
$ws.Range("M147").Value = 2.1048
$ws.Range("N147").Value = 2.7462
$ws.Range("O147").Value = 1784937.706
$ws.Range("P147").Value = 5565566.687
$ws.Range("Q147").Value = "Rangitikei District"
$ws.Range("R147").Value = "Rangitīkei-Turakina"
$ws.Range("S147").Value = "Turakina"
$ws.Range("T147").Value = "Tura_1c"
$ws.Range("U147").Value = "mg/L"
# Row 148
$ws.Range("A148").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B148").Value = "E coli (>260)"
$ws.Range("C148").Value = "E"
$ws.Range("D148").Value = "2019 - 2023"
$ws.Range("E148").Value = "Impact"
$ws.Range("F148").Value = 2000
$ws.Range("G148").Value = 6768.49687927857
$ws.Range("H148").Value = 73000
$ws.Range("I148").Value = 36489
$ws.Range("J148").Value = 63.1578947368421
$ws.Range("K148").Value = 77.1929824561404
$ws.Range("L148").Value = 330
$ws.Range("M148").Value = 9700
$ws.Range("N148").Value = 19295.56279
$ws.Range("O148").Value = 1784937.706
$ws.Range("P148").Value = 5565566.687
$ws.Range("Q148").Value = "Rangitikei District"
$ws.Range("R148").Value = "Rangitīkei-Turakina"
$ws.Range("S148").Value = "Turakina"
$ws.Range("T148").Value = "Tura_1c"
$ws.Range("U148").Value = "% exceedances over 260/100 mL"
# Row 149
$ws.Range("A149").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B149").Value = "E coli (>540)"
$ws.Range("C149").Value = "E"
$ws.Range("D149").Value = "2019 - 2023"
$ws.Range("E149").Value = "Impact"
$ws.Range("F149").Value = 2000
$ws.Range("G149").Value = 6768.49687927857
$ws.Range("H149").Value = 73000
$ws.Range("I149").Value = 36489
$ws.Range("J149").Value = 63.1578947368421
$ws.Range("K149").Value = 77.1929824561404
$ws.Range("L149").Value = 330
$ws.Range("M149").Value = 9700
$ws.Range("N149").Value = 19295.56279
$ws.Range("O149").Value = 1784937.706
$ws.Range("P149").Value = 5565566.687
$ws.Range("Q149").Value = "Rangitikei District"
$ws.Range("R149").Value = "Rangitīkei-Turakina"
$ws.Range("S149").Value = "Turakina"
$ws.Range("T149").Value = "Tura_1c"
$ws.Range("U149").Value = "% exceedances over 540/100 mL"
# Row 150
$ws.Range("A150").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B150").Value = "E coli (Median)"
$ws.Range("C150").Value = "E"
$ws.Range("D150").Value = "2019 - 2023"
$ws.Range("E150").Value = "Impact"
$ws.Range("F150").Value = 2000
$ws.Range("G150").Value = 6768.49687927857
$ws.Range("H150").Value = 73000
$ws.Range("I150").Value = 36489
$ws.Range("J150").Value = 63.1578947368421
$ws.Range("K150").Value = 77.1929824561404
$ws.Range("L150").Value = 330
$ws.Range("M150").Value = 9700
$ws.Range("N150").Value = 19295.56279
$ws.Range("O150").Value = 1784937.706
$ws.Range("P150").Value = 5565566.687
$ws.Range("Q150").Value = "Rangitikei District"
$ws.Range("R150").Value = "Rangitīkei-Turakina"
$ws.Range("S150").Value = "Turakina"
$ws.Range("T150").Value = "Tura_1c"
$ws.Range("U150").Value = "E. coli/100 mL"
# Row 151
$ws.Range("A151").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B151").Value = "E coli (95th Percentile)"
$ws.Range("C151").Value = "E"
$ws.Range("D151").Value = "2019 - 2023"
$ws.Range("E151").Value = "Impact"
$ws.Range("F151").Value = 2000
$ws.Range("G151").Value = 6768.49687927857
$ws.Range("H151").Value = 73000
$ws.Range("I151").Value = 36489
$ws.Range("J151").Value = 63.1578947368421
$ws.Range("K151").Value = 77.1929824561404
$ws.Range("L151").Value = 330
$ws.Range("M151").Value = 9700
$ws.Range("N151").Value = 19295.56279
$ws.Range("O151").Value = 1784937.706
$ws.Range("P151").Value = 5565566.687
$ws.Range("Q151").Value = "Rangitikei District"
$ws.Range("R151").Value = "Rangitīkei-Turakina"
$ws.Range("S151").Value = "Turakina"
$ws.Range("T151").Value = "Tura_1c"
$ws.Range("U151").Value = "E. coli/100 mL"
# Row 152
$ws.Range("A152").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B152").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C152").Value = "D"
$ws.Range("D152").Value = "2019 - 2023"
$ws.Range("E152").Value = "Impact"
$ws.Range("F152").Value = 5.32021
$ws.Range("G152").Value = 8.294741192802199
$ws.Range("H152").Value = 37.6
$ws.Range("I152").Value = 29.97029
$ws.Range("L152").Value = 1.6925
$ws.Range("M152").Value = 16.337
$ws.Range("N152").Value = 26.16606
$ws.Range("O152").Value = 1784937.706
$ws.Range("P152").Value = 5565566.687
$ws.Range("Q152").Value = "Rangitikei District"
$ws.Range("R152").Value = "Rangitīkei-Turakina"
$ws.Range("S152").Value = "Turakina"
$ws.Range("T152").Value = "Tura_1c"
$ws.Range("U152").Value = "mg NH4-N/L"
# Row 153
$ws.Range("A153").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B153").Value = "Ammoniacal-N (Median)"
$ws.Range("C153").Value = "D"
$ws.Range("D153").Value = "2019 - 2023"
$ws.Range("E153").Value = "Impact"
$ws.Range("F153").Value = 5.32021
$ws.Range("G153").Value = 8.294741192802199
$ws.Range("H153").Value = 37.6
$ws.Range("I153").Value = 29.97029
$ws.Range("L153").Value = 1.6925
$ws.Range("M153").Value = 16.337
$ws.Range("N153").Value = 26.16606
$ws.Range("O153").Value = 1784937.706
$ws.Range("P153").Value = 5565566.687
$ws.Range("Q153").Value = "Rangitikei District"
$ws.Range("R153").Value = "Rangitīkei-Turakina"
$ws.Range("S153").Value = "Turakina"
$ws.Range("T153").Value = "Tura_1c"
$ws.Range("U153").Value = "mg NH4-N/L"
# Row 154
$ws.Range("A154").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B154").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C154").Value = "B"
$ws.Range("D154").Value = "2019 - 2023"
$ws.Range("E154").Value = "Impact"
$ws.Range("F154").Value = 0.64
$ws.Range("G154").Value = 0.786210526315789
$ws.Range("H154").Value = 2.6
$ws.Range("I154").Value = 2.036
$ws.Range("L154").Value = 0.532
$ws.Range("M154").Value = 1.3024
$ws.Range("N154").Value = 1.6194
$ws.Range("O154").Value = 1784937.706
$ws.Range("P154").Value = 5565566.687
$ws.Range("Q154").Value = "Rangitikei District"
$ws.Range("R154").Value = "Rangitīkei-Turakina"
$ws.Range("S154").Value = "Turakina"
$ws.Range("T154").Value = "Tura_1c"
$ws.Range("U154").Value = "mg NO3-N/L"
# Row 155
$ws.Range("A155").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B155").Value = "Nitrate-N (Median)"
$ws.Range("C155").Value = "A"
$ws.Range("D155").Value = "2019 - 2023"
$ws.Range("E155").Value = "Impact"
$ws.Range("F155").Value = 0.64
$ws.Range("G155").Value = 0.786210526315789
$ws.Range("H155").Value = 2.6
$ws.Range("I155").Value = 2.036
$ws.Range("L155").Value = 0.532
$ws.Range("M155").Value = 1.3024
$ws.Range("N155").Value = 1.6194
$ws.Range("O155").Value = 1784937.706
$ws.Range("P155").Value = 5565566.687
$ws.Range("Q155").Value = "Rangitikei District"
$ws.Range("R155").Value = "Rangitīkei-Turakina"
$ws.Range("S155").Value = "Turakina"
$ws.Range("T155").Value = "Tura_1c"
$ws.Range("U155").Value = "mg NO3-N/L"
# Row 156
$ws.Range("A156").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B156").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("D156").Value = "2019 - 2023"
$ws.Range("E156").Value = "Impact"
$ws.Range("F156").Value = 5.99
$ws.Range("G156").Value = 7.32724561403509
$ws.Range("H156").Value = 22.71
$ws.Range("I156").Value = 20.887
$ws.Range("L156").Value = 1.782
$ws.Range("M156").Value = 13.1003
$ws.Range("N156").Value = 19.2262
$ws.Range("O156").Value = 1784937.706
$ws.Range("P156").Value = 5565566.687
$ws.Range("Q156").Value = "Rangitikei District"
$ws.Range("R156").Value = "Rangitīkei-Turakina"
$ws.Range("S156").Value = "Turakina"
$ws.Range("T156").Value = "Tura_1c"
$ws.Range("U156").Value = "g/m3"
# Row 157
$ws.Range("A157").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B157").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("D157").Value = "2019 - 2023"
$ws.Range("E157").Value = "Impact"
$ws.Range("F157").Value = 5.99
$ws.Range("G157").Value = 7.32724561403509
$ws.Range("H157").Value = 22.71
$ws.Range("I157").Value = 20.887
$ws.Range("L157").Value = 1.782
$ws.Range("M157").Value = 13.1003
$ws.Range("N157").Value = 19.2262
$ws.Range("O157").Value = 1784937.706
$ws.Range("P157").Value = 5565566.687
$ws.Range("Q157").Value = "Rangitikei District"
$ws.Range("R157").Value = "Rangitīkei-Turakina"
$ws.Range("S157").Value = "Turakina"
$ws.Range("T157").Value = "Tura_1c"
$ws.Range("U157").Value = "g/m3"
# Row 158
$ws.Range("A158").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B158").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("D158").Value = "2019 - 2023"
$ws.Range("E158").Value = "Impact"
$ws.Range("F158").Value = 12.1
$ws.Range("G158").Value = 12.629298245614
$ws.Range("H158").Value = 24.7
$ws.Range("I158").Value = 23.49
$ws.Range("L158").Value = 11.3
$ws.Range("M158").Value = 19.224
$ws.Range("N158").Value = 23.082
$ws.Range("O158").Value = 1784937.706
$ws.Range("P158").Value = 5565566.687
$ws.Range("Q158").Value = "Rangitikei District"
$ws.Range("R158").Value = "Rangitīkei-Turakina"
$ws.Range("S158").Value = "Turakina"
$ws.Range("T158").Value = "Tura_1c"
$ws.Range("U158").Value = "g/m3"
# Row 159
$ws.Range("A159").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B159").Value = "Total Nitrogen (Median)"
$ws.Range("D159").Value = "2019 - 2023"
$ws.Range("E159").Value = "Impact"
$ws.Range("F159").Value = 12.1
$ws.Range("G159").Value = 12.629298245614
$ws.Range("H159").Value = 24.7
$ws.Range("I159").Value = 23.49
$ws.Range("L159").Value = 11.3
$ws.Range("M159").Value = 19.224
$ws.Range("N159").Value = 23.082
$ws.Range("O159").Value = 1784937.706
$ws.Range("P159").Value = 5565566.687
$ws.Range("Q159").Value = "Rangitikei District"
$ws.Range("R159").Value = "Rangitīkei-Turakina"
$ws.Range("S159").Value = "Turakina"
$ws.Range("T159").Value = "Tura_1c"
$ws.Range("U159").Value = "g/m3"
# Row 160
$ws.Range("A160").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B160").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("D160").Value = "2019 - 2023"
$ws.Range("E160").Value = "Impact"
$ws.Range("F160").Value = 2.14
$ws.Range("G160").Value = 2.13522807017544
$ws.Range("H160").Value = 5.06
$ws.Range("I160").Value = 4.143
$ws.Range("L160").Value = 2.23
$ws.Range("M160").Value = 3.1981
$ws.Range("N160").Value = 3.8334
$ws.Range("O160").Value = 1784937.706
$ws.Range("P160").Value = 5565566.687
$ws.Range("Q160").Value = "Rangitikei District"
$ws.Range("R160").Value = "Rangitīkei-Turakina"
$ws.Range("S160").Value = "Turakina"
$ws.Range("T160").Value = "Tura_1c"
$ws.Range("U160").Value = "g/m3"
# Row 161
$ws.Range("A161").Value = "Unnamed Trib of Waipu at ds Ratana STP"
$ws.Range("B161").Value = "Total Phosphorus (Median)"
$ws.Range("D161").Value = "2019 - 2023"
$ws.Range("E161").Value = "Impact"
$ws.Range("F161").Value = 2.14
$ws.Range("G161").Value = 2.13522807017544
$ws.Range("H161").Value = 5.06
$ws.Range("I161").Value = 4.143
$ws.Range("L161").Value = 2.23
$ws.Range("M161").Value = 3.1981
$ws.Range("N161").Value = 3.8334
$ws.Range("O161").Value = 1784937.706
$ws.Range("P161").Value = 5565566.687
$ws.Range("Q161").Value = "Rangitikei District"
$ws.Range("R161").Value = "Rangitīkei-Turakina"
$ws.Range("S161").Value = "Turakina"
$ws.Range("T161").Value = "Tura_1c"
$ws.Range("U161").Value = "g/m3"
